# Applies the cryptos.xlsx symbol-list update described in the commit:
# "Updated symbol list on Tue Feb 14 04:58:50 UTC 2023 with GitHub Actions"
#
# GateToken moved up from row 17 to row 6 (re-ranked), shifting the rows
# for FTXToken..LEO down by one, and fresh Price/Volume(1h) figures were
# pulled for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Price/Volume columns (D, E) hold numeric-looking text
# (e.g. "292.07", "-6.50%") in the source workbook, so a leading apostrophe
# is used to force Excel to store them as text instead of reinterpreting
# them as a number/percentage.
$updates = [ordered]@{
    "D2" = "'292.07"
    "E2" = "'-6.50%"
    "D3" = "'40.27"
    "E3" = "'0.85%"
    "D4" = "'5.004"
    "E4" = "'-3.54%"
    "D5" = "'0.07317"
    "E5" = "'-3.46%"
    "B6" = "GateToken"
    "C6" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D6" = "'4.295"
    "E6" = "'-0.69%"
    "B7" = "FTXToken"
    "C7" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D7" = "'1.521"
    "E7" = "'-8.35%"
    "B8" = "MXToken"
    "C8" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8" = "'0.9275"
    "E8" = "'0.17%"
    "B9" = "BTSEToken"
    "C9" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D9" = "'2.378"
    "E9" = "'-1.90%"
    "B10" = "LiechtensteinCryptoassetsExchange"
    "C10" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D10" = "'0.1192"
    "E10" = "'-0.65%"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D11" = "'0.1737"
    "E11" = "'-4.99%"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "'0.04331"
    "E12" = "'3.96%"
    "B13" = "MandalaExchangeToken"
    "C13" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D13" = "'0.08632"
    "E13" = "'-4.48%"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "'0.1054"
    "E14" = "'0.07%"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "'0.001271"
    "E15" = "'-1.28%"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D16" = "'0.005958"
    "E16" = "'2.03%"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D17" = "'3.339"
    "E17" = "'-0.31%"
    "E18" = "'-2.01%"
    "D19" = "'7.976"
    "E19" = "'5.03%"
    "E20" = "'2.90%"
    "E21" = "'-0.59%"
    "D22" = "'0.03932"
    "E22" = "'-2.13%"
    "E23" = "'-0.61%"
    "E24" = "'-7.16%"
    "E25" = "'0.89%"
    "E26" = "'-95.05%"
    "D38" = "'0.02286"
    "E38" = "'-5.62%"
    "D39" = "'0.04980"
    "E39" = "'-3.46%"
    "D40" = "'0.005357"
    "E40" = "'62.29%"
    "D41" = "'0.007701"
    "E41" = "'-0.60%"
    "D42" = "'0.1283"
    "E42" = "'-1.21%"
    "D43" = "'0.007329"
    "E43" = "'-3.66%"
    "D44" = "'0.007891"
    "E44" = "'-3.77%"
    "D45" = "'0.3170"
    "E45" = "'2.15%"
    "D46" = "'0.00006322"
    "E46" = "'-4.02%"
    "E47" = "'0.11%"
    "D48" = "'0.02043"
    "E48" = "'-92.71%"
    "E49" = "'0.11%"
    "E50" = "'0.11%"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
